$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update summary header figures (VALOR MORA total, worker count, period count)
$ws.Range("E11").Value = 58673
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 2

# Update existing first data row's "Salario Basico" (row 16 - Oscar Evelio Loaiza Noguera)
$ws.Range("G16").Value = 1300000

# Remove the intermediate data rows for the old workers (Keiner, Javier, Danis - all but
# the very last of their period rows). The surviving last row (30) slides up to become
# row 17 and keeps its distinctive "closing row" border/style set.
$ws.Rows("17:29").Delete()

# Replace that now-row-17 with the new worker's data (Felipe Andres Cardona Ariza)
$ws.Range("C17").Value = "1007855507"
$ws.Range("D17").Value = "FELIPE ANDRES CARDONA ARIZA"
$ws.Range("E17").Value = "2508"
